$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells being updated with text-like values that could otherwise be
# auto-coerced to numbers by Excel (e.g. "585.46"). Force text format,
# write the value, then restore the default "Normal" style so the
# cell keeps looking like the rest of the sheet (no explicit style).
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E8", "E10", "D11", "E11", "E12", "D13", "E13", "D14", "E14", "E16", "D17", "E17", "D18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "E27", "E28", "D29", "E29", "D30", "E30", "E31", "D32", "E32", "E33", "D34", "E34", "E35", "D36", "E36", "E37", "E38", "D39", "E39", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "E49", "D50", "E50", "E51")
foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '64.156.39'
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '3.521.68'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '585.46'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("D6").Value = '134.58'
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("D7").Value = '3.521.87'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").Value = '7.13'
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("E12").Value = '  -2.13%  '
$ws.Range("D13").Value = '4.122.60'
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").Value = '27.49'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("D17").Value = '3.520.76'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '64.193.07'
$ws.Range("D19").Value = '9.78'
$ws.Range("E19").Value = '  -2.77%  '
$ws.Range("D20").Value = '13.90'
$ws.Range("E20").Value = '  -2.89%  '
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("D22").Value = '382.84'
$ws.Range("E22").Value = '  -2.10%  '
$ws.Range("D23").Value = '0.571'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("D24").Value = '3.665.39'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").Value = '74.05'
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  -1.29%  '
$ws.Range("E28").Value = '  +3.60%  '
$ws.Range("D29").Value = '1.58'
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("D30").Value = '7.46'
$ws.Range("E30").Value = '  -2.06%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").Value = '8.44'
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("D34").Value = '3.535.79'
$ws.Range("E34").Value = '  +0.31%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '23.60'
$ws.Range("E36").Value = '  -2.05%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("E38").Value = '  +3.96%  '
$ws.Range("D39").Value = '6.95'
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").Value = '160.52'
$ws.Range("E41").Value = '  -5.32%  '
$ws.Range("D42").Value = '0.0786'
$ws.Range("E42").Value = '  -2.33%  '
$ws.Range("D43").Value = '26.73'
$ws.Range("E43").Value = '  +2.44%  '
$ws.Range("D44").Value = '0.812'
$ws.Range("E44").Value = '  -1.04%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").Value = '  -2.73%  '
$ws.Range("D47").Value = '41.73'
$ws.Range("E47").Value = '  -2.74%  '
$ws.Range("D48").Value = '4.40'
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("E49").Value = '  -2.54%  '
$ws.Range("D50").Value = '2.484.70'
$ws.Range("E50").Value = '  +1.18%  '
$ws.Range("E51").Value = '  -1.03%  '

foreach ($c in $cells) {
    $ws.Range($c).Style = "Normal"
}
